$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column O entirely (U.S. Trump column) - this shifts nothing else,
# we just delete the last used column.
$ws.Columns.Item(15).Delete()

# Update header row (row 1). A1 and B1 are unchanged.
$ws.Range("C1").Value = "`$ bold('Europe')"
$ws.Range("D1").Value = "Saudi Arabia"
$ws.Range("E1").Value = "Italy"
$ws.Range("F1").Value = "Spain"
$ws.Range("G1").Value = "Germany"
$ws.Range("H1").Value = "United Kingdom"
$ws.Range("I1").Value = "Russia"
$ws.Range("J1").Value = "France"
$ws.Range("K1").Value = "Poland"
$ws.Range("L1").Value = "USA"
$ws.Range("M1").Value = "Switzerland"
$ws.Range("N1").Value = "Japan"

# Update data rows 2-6, columns C-N (A and B unchanged).
$ws.Range("C2").Value = 0.111078823224324
$ws.Range("D2").Value = 0.486987142752865
$ws.Range("E2").Value = 0.394246389814818
$ws.Range("F2").Value = 0.226787965621744
$ws.Range("G2").Value = 0.0762722737213013
$ws.Range("H2").Value = 0.0480759539306938
$ws.Range("I2").Value = 0.0410217438167706
$ws.Range("J2").Value = 0.0178306762563666
$ws.Range("K2").Value = -0.0854313280129163
$ws.Range("L2").Value = -0.104086448159795
$ws.Range("M2").Value = -0.114323136373101
$ws.Range("N2").Value = -0.209382614568853

$ws.Range("C3").Value = 0.563770192004184
$ws.Range("D3").Value = 0.637895863088784
$ws.Range("E3").Value = 0.651321187680723
$ws.Range("F3").Value = 0.580138863192046
$ws.Range("G3").Value = 0.553113452044712
$ws.Range("H3").Value = 0.552228264162119
$ws.Range("I3").Value = 0.496693987290596
$ws.Range("J3").Value = 0.546518787466243
$ws.Range("K3").Value = 0.486594535915121
$ws.Range("L3").Value = 0.484917310839545
$ws.Range("M3").Value = 0.525215324624544
$ws.Range("N3").Value = 0.382887811020303

$ws.Range("C4").Value = 0.20593904313023
$ws.Range("D4").Value = 0.142379895004253
$ws.Range("E4").Value = 0.156515577986039
$ws.Range("F4").Value = 0.194393323940455
$ws.Range("G4").Value = 0.213519037471354
$ws.Range("H4").Value = 0.209554744465343
$ws.Range("I4").Value = 0.182723223867721
$ws.Range("J4").Value = 0.220108318263721
$ws.Range("K4").Value = 0.233463656969949
$ws.Range("L4").Value = 0.230976796140146
$ws.Range("M4").Value = 0.283725963176114
$ws.Range("N4").Value = 0.18122035142742

$ws.Range("C5").Value = 0.357831148873954
$ws.Range("D5").Value = 0.495515968084531
$ws.Range("E5").Value = 0.494805609694684
$ws.Range("F5").Value = 0.385745539251591
$ws.Range("G5").Value = 0.339594414573358
$ws.Range("H5").Value = 0.342673519696776
$ws.Range("I5").Value = 0.313970763422875
$ws.Range("J5").Value = 0.326410469202522
$ws.Range("K5").Value = 0.253130878945172
$ws.Range("L5").Value = 0.253940514699399
$ws.Range("M5").Value = 0.241489361448431
$ws.Range("N5").Value = 0.201667459592883

$ws.Range("C6").Value = 0.720831013859595
$ws.Range("D6").Value = 0.798467109724214
$ws.Range("E6").Value = 0.785932597541719
$ws.Range("F6").Value = 0.737279969748424
$ws.Range("G6").Value = 0.711660103572645
$ws.Range("H6").Value = 0.712225124939768
$ws.Range("I6").Value = 0.711487796713406
$ws.Range("J6").Value = 0.70760932096738
$ws.Range("K6").Value = 0.671297716035589
$ws.Range("L6").Value = 0.675107606459803
$ws.Range("M6").Value = 0.647088859608711
$ws.Range("N6").Value = 0.673290897843638
